$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the existing header style (bold, centered, bordered) from AC1 onto the
# new header cells before setting their text, so AD1:AF1 match the rest of
# the header row (s="1").
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Header row: add Wins / Losses / Ties in columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team's W/L/T record for every player row (2-54)
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 30).Value = 79   # AD = Wins
    $ws.Cells.Item($r, 31).Value = 83   # AE = Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = Ties
}
